$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string rich-text cells: update the "volume/number" banner and the report week banner.
# Both strings are built of runs sharing identical formatting, so writing the
# full text back in one go renders identically to editing just the one run.
$ws.Range("C1").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G14").PasteSpecial(-4122)

$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)

$ws.Range("D20").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Value = 0
$ws.Range("K14").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("C22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G22").PasteSpecial(-4122)

$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H22").PasteSpecial(-4122)

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4122)

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E26").PasteSpecial(-4122)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("D28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("D29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
# --- Simple numeric value updates (same type/style, value only) ---
$ws.Range("M14").Value = -66.666666666666
$ws.Range("F15").Value = 6
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 13
$ws.Range("K15").Value = 18.181818181818
$ws.Range("L15").Value = 333.333333333333
$ws.Range("M15").Value = 18.181818181818
$ws.Range("N15").Value = -51.851851851851
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 14.285714285714
$ws.Range("I16").Value = 73
$ws.Range("J16").Value = 64
$ws.Range("K16").Value = 14.0625
$ws.Range("L16").Value = -15.116279069767
$ws.Range("M16").Value = -40.163934426229
$ws.Range("N16").Value = -90.161725067385
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 39
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 21.875
$ws.Range("I17").Value = 140
$ws.Range("J17").Value = 148
$ws.Range("K17").Value = -5.405405405405
$ws.Range("L17").Value = 3.703703703703
$ws.Range("M17").Value = -10.828025477707
$ws.Range("N17").Value = -66.101694915254
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 30.769230769230
$ws.Range("I18").Value = 76
$ws.Range("J18").Value = 76
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 26.666666666666
$ws.Range("M18").Value = -50.649350649350
$ws.Range("N18").Value = -82.882882882882
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 38.461538461538
$ws.Range("I19").Value = 155
$ws.Range("J19").Value = 127
$ws.Range("K19").Value = 22.047244094488
$ws.Range("L19").Value = 28.099173553719
$ws.Range("M19").Value = 7.638888888888
$ws.Range("N19").Value = -44.244604316546
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 14.285714285714
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 35
$ws.Range("K20").Value = 34.285714285714
$ws.Range("L20").Value = 38.235294117647
$ws.Range("M20").Value = 20.512820512820
$ws.Range("N20").Value = -78.139534883720
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 6.896551724137
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = 29.787234042553
$ws.Range("I21").Value = 506
$ws.Range("J21").Value = 463
$ws.Range("K21").Value = 9.287257019438
$ws.Range("L21").Value = 14.479638009049
$ws.Range("M21").Value = -20.063191153238
$ws.Range("N21").Value = -76.355140186915
$ws.Range("I22").Value = 6
$ws.Range("K22").Value = -40
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -53.846153846153
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = -16.666666666666
$ws.Range("F23").Value = 22
$ws.Range("H23").Value = 22.222222222222
$ws.Range("I23").Value = 87
$ws.Range("J23").Value = 73
$ws.Range("K23").Value = 19.178082191780
$ws.Range("L23").Value = -1.136363636363
$ws.Range("M23").Value = 19.178082191780
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -27.272727272727
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 117
$ws.Range("H24").Value = -16.239316239316
$ws.Range("I24").Value = 551
$ws.Range("J24").Value = 558
$ws.Range("K24").Value = -1.254480286738
$ws.Range("L24").Value = 32.134292565947
$ws.Range("M24").Value = 69.018404907975
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = -5.263157894736
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = -5.882352941176
$ws.Range("I25").Value = 237
$ws.Range("J25").Value = 245
$ws.Range("K25").Value = -3.265306122448
$ws.Range("L25").Value = 54.901960784313
$ws.Range("M25").Value = -23.300970873786
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 7
$ws.Range("H26").Value = 75
$ws.Range("I26").Value = 20
$ws.Range("K26").Value = 17.647058823529
$ws.Range("L26").Value = 150
$ws.Range("G27").Value = 4
$ws.Range("L27").Value = -50
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 50
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = -45.454545454545
$ws.Range("L28").Value = -64.705882352941
$ws.Range("M28").Value = -72.727272727272
$ws.Range("N28").Value = -94.594594594594
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("J29").Value = 8
$ws.Range("K29").Value = -37.5
$ws.Range("L29").Value = -66.666666666666
$ws.Range("M29").Value = -68.75
$ws.Range("N29").Value = -94.565217391304
